$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume figures.
# For values in column D that look like plain numbers, a leading
# apostrophe is used to force Excel to keep them as text (matching
# the original inline-string cell type) instead of converting them
# into numeric cells.

$ws.Range("D2").Value = "58.998.65"
$ws.Range("E2").Value = "  -3.19%  "
$ws.Range("D3").Value = "3.230.86"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'540.84"
$ws.Range("E5").Value = "  -4.50%  "
$ws.Range("D6").Value = "'136.44"
$ws.Range("D8").Value = "3.228.66"
$ws.Range("E8").Value = "  -3.90%  "
$ws.Range("E9").Value = "  -4.11%  "
$ws.Range("D10").Value = "'7.64"
$ws.Range("E10").Value = "  -4.08%  "
$ws.Range("E11").Value = "  -5.66%  "
$ws.Range("D12").Value = "'0.395"
$ws.Range("E12").Value = "  -4.31%  "
$ws.Range("D13").Value = "3.782.81"
$ws.Range("E13").Value = "  -3.89%  "
$ws.Range("E14").Value = "  -0.87%  "
$ws.Range("D15").Value = "'26.02"
$ws.Range("E15").Value = "  -7.08%  "
$ws.Range("D16").Value = "3.229.34"
$ws.Range("E16").Value = "  -3.75%  "
$ws.Range("E17").Value = "  -5.90%  "
$ws.Range("D18").Value = "59.013.35"
$ws.Range("E18").Value = "  -3.28%  "
$ws.Range("D19").Value = "'5.92"
$ws.Range("E19").Value = "  -6.77%  "
$ws.Range("D20").Value = "'13.38"
$ws.Range("E20").Value = "  -5.85%  "
$ws.Range("D22").Value = "'363.25"
$ws.Range("E22").Value = "  -3.12%  "
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").Value = "'70.62"
$ws.Range("E24").Value = "  -6.17%  "
$ws.Range("D25").Value = "'0.522"
$ws.Range("E25").Value = "  -6.69%  "
$ws.Range("D26").Value = "3.362.57"
$ws.Range("E26").Value = "  -3.93%  "
$ws.Range("D27").Value = "'0.172"
$ws.Range("E27").Value = "  -1.81%  "
$ws.Range("D28").Value = "'0.0$([char]8323)0973"
$ws.Range("E28").Value = "  -10.14%  "
$ws.Range("D29").Value = "'1.01"
$ws.Range("E29").Value = "  +0.85%  "
$ws.Range("D30").Value = "'7.19"
$ws.Range("E30").Value = "  -2.75%  "
$ws.Range("E32").Value = "  -6.56%  "
$ws.Range("D33").Value = "'7.13"
$ws.Range("E33").Value = "  -7.36%  "
$ws.Range("D34").Value = "'21.94"
$ws.Range("E34").Value = "  -3.98%  "
$ws.Range("E35").Value = "  -5.82%  "
$ws.Range("D36").Value = "'4.96"
$ws.Range("E36").Value = "  -7.56%  "
$ws.Range("D37").Value = "'161.44"
$ws.Range("E37").Value = "  -5.00%  "
$ws.Range("D38").Value = "'6.44"
$ws.Range("E38").Value = "  -5.00%  "
$ws.Range("D39").Value = "'1.45"
$ws.Range("E39").Value = "  -5.99%  "
$ws.Range("D40").Value = "'26.37"
$ws.Range("E40").Value = "  -9.18%  "
$ws.Range("D41").Value = "'0.0711"
$ws.Range("E41").Value = "  -4.64%  "
$ws.Range("D42").Value = "3.258.38"
$ws.Range("E42").Value = "  -4.03%  "
$ws.Range("D43").Value = "'41.15"
$ws.Range("E43").Value = "  -2.72%  "
$ws.Range("E44").Value = "  -5.72%  "
$ws.Range("D45").Value = "'1.11"
$ws.Range("E45").Value = "  -3.25%  "
$ws.Range("E46").Value = "  -5.62%  "
$ws.Range("E47").Value = "  -6.01%  "
$ws.Range("E48").Value = "  -0.04%  "
$ws.Range("D49").Value = "2.304.31"
$ws.Range("E49").Value = "  -7.34%  "
$ws.Range("D50").Value = "'6.33"
$ws.Range("E50").Value = "  -4.92%  "
$ws.Range("D51").Value = "'20.98"
$ws.Range("E51").Value = "  -6.80%  "
